# Updates market-price-derived columns (H:N) on several leve rows across
# the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to match a refreshed price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1808.5974
$ws.Range("I17").Value = 427.27274
$ws.Range("J17").Value = 2038.8182
$ws.Range("K17").Value = 1281.81822
$ws.Range("L17").Value = 6116.4546
$ws.Range("M17").Value = -1113.81822
$ws.Range("N17").Value = -6452.4546
$ws.Range("H62").Value = 3613.054
$ws.Range("I62").Value = 2784
$ws.Range("K62").Value = 2784
$ws.Range("M62").Value = -2160
$ws.Range("H64").Value = 7062.625
$ws.Range("I64").Value = 6666.6665
$ws.Range("J64").Value = 7300.2
$ws.Range("K64").Value = 6666.6665
$ws.Range("L64").Value = 7300.2
$ws.Range("M64").Value = -6418.6665
$ws.Range("N64").Value = -7796.2
$ws.Range("H65").Value = 3613.054
$ws.Range("I65").Value = 2784
$ws.Range("K65").Value = 13920
$ws.Range("M65").Value = -10800
$ws.Range("H67").Value = 7062.625
$ws.Range("I67").Value = 6666.6665
$ws.Range("J67").Value = 7300.2
$ws.Range("K67").Value = 6666.6665
$ws.Range("L67").Value = 7300.2
$ws.Range("M67").Value = -5808.6665
$ws.Range("N67").Value = -9016.200000000001
$ws.Range("H70").Value = 5512.3335
$ws.Range("H73").Value = 5512.3335
$ws.Range("H107").Value = 2165.625
$ws.Range("I107").Value = 1772.64
$ws.Range("J107").Value = 3569.1428
$ws.Range("K107").Value = 1772.64
$ws.Range("L107").Value = 3569.1428
$ws.Range("M107").Value = 147.3599999999999
$ws.Range("N107").Value = -7409.1428
$ws.Range("H127").Value = 567.6
$ws.Range("I127").Value = 567.6
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 1702.8
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = 3257.2
$ws.Range("M127").ClearContents()
$ws.Range("H137").Value = 3432.3076
$ws.Range("I137").Value = 2758
$ws.Range("J137").Value = 3732
$ws.Range("K137").Value = 8274
$ws.Range("L137").Value = 11196
$ws.Range("M137").Value = -5724
$ws.Range("N137").Value = -16296

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2742.8164
$ws.Range("I32").Value = 2255.442
$ws.Range("J32").Value = 6235.6665
$ws.Range("K32").Value = 2255.442
$ws.Range("L32").Value = 6235.6665
$ws.Range("M32").Value = -1968.442
$ws.Range("N32").Value = -6809.6665
$ws.Range("H45").Value = 52634292
$ws.Range("I45").Value = 71429910
$ws.Range("K45").Value = 71429910
$ws.Range("M45").Value = -71429533
$ws.Range("H63").Value = 2833
$ws.Range("H66").Value = 2833
$ws.Range("H74").Value = 22224426
$ws.Range("I74").Value = 23811742
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 23811742
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -23810868
$ws.Range("N74").Value = -3748
$ws.Range("H77").Value = 22224426
$ws.Range("I77").Value = 23811742
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 119058710
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -119054342
$ws.Range("N77").Value = -18736
$ws.Range("H110").Value = 4519.3335
$ws.Range("I110").Value = 1884.6666
$ws.Range("J110").Value = 7154
$ws.Range("K110").Value = 1884.6666
$ws.Range("L110").Value = 7154
$ws.Range("M110").Value = 160.3334
$ws.Range("N110").Value = -11244
$ws.Range("H122").Value = 2029.7693
$ws.Range("I122").Value = 996.5
$ws.Range("J122").Value = 2915.4285
$ws.Range("K122").Value = 2989.5
$ws.Range("L122").Value = 8746.2855
$ws.Range("M122").Value = -539.5
$ws.Range("N122").Value = -13646.2855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 60801.8
$ws.Range("I105").Value = 111505
$ws.Range("K105").Value = 111505
$ws.Range("M105").Value = -109758
$ws.Range("H107").Value = 633.625
$ws.Range("I107").Value = 409.78946
$ws.Range("K107").Value = 409.78946
$ws.Range("M107").Value = 1510.21054
$ws.Range("H134").Value = 2952.3914
$ws.Range("I134").Value = 1873.6316
$ws.Range("J134").Value = 8076.5
$ws.Range("K134").Value = 5620.8948
$ws.Range("L134").Value = 24229.5
$ws.Range("M134").Value = -3085.8948
$ws.Range("N134").Value = -29299.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 35086.516
$ws.Range("I31").Value = 3042.5833
$ws.Range("J31").Value = 120537
$ws.Range("K31").Value = 3042.5833
$ws.Range("L31").Value = 120537
$ws.Range("M31").Value = -2747.5833
$ws.Range("N31").Value = -121127
$ws.Range("H34").Value = 35086.516
$ws.Range("I34").Value = 3042.5833
$ws.Range("J34").Value = 120537
$ws.Range("K34").Value = 3042.5833
$ws.Range("L34").Value = 120537
$ws.Range("M34").Value = -2840.5833
$ws.Range("N34").Value = -120941
$ws.Range("H58").Value = 3140.1482
$ws.Range("I58").Value = 1198.05
$ws.Range("K58").Value = 1198.05
$ws.Range("M58").Value = -995.05
$ws.Range("H64").Value = 25000
$ws.Range("J64").Value = 25000
$ws.Range("L64").Value = 25000
$ws.Range("N64").Value = -25496
$ws.Range("H67").Value = 25000
$ws.Range("J67").Value = 25000
$ws.Range("L67").Value = 25000
$ws.Range("N67").Value = -26716
$ws.Range("H74").Value = 53888
$ws.Range("I74").Value = 7777
$ws.Range("K74").Value = 7777
$ws.Range("M74").Value = -6903
$ws.Range("H77").Value = 53888
$ws.Range("I77").Value = 7777
$ws.Range("K77").Value = 23331
$ws.Range("M77").Value = -18963
$ws.Range("H134").Value = 2442.6897
$ws.Range("I134").Value = 1651.9524
$ws.Range("J134").Value = 4518.375
$ws.Range("K134").Value = 4955.857199999999
$ws.Range("L134").Value = 13555.125
$ws.Range("M134").Value = -2420.857199999999
$ws.Range("N134").Value = -18625.125
$ws.Range("H136").Value = 3140.1482
$ws.Range("I136").Value = 1198.05
$ws.Range("K136").Value = 3594.15
$ws.Range("M136").Value = -1044.15

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 5569.7407
$ws.Range("J5").Value = 9964.714
$ws.Range("L5").Value = 29894.142
$ws.Range("N5").Value = -30118.142
$ws.Range("H33").Value = 72.72221999999999
$ws.Range("I33").Value = 74.69231000000001
$ws.Range("J33").Value = 67.59999999999999
$ws.Range("K33").Value = 448.15386
$ws.Range("L33").Value = 405.6
$ws.Range("M33").Value = -165.15386
$ws.Range("N33").Value = -971.5999999999999
$ws.Range("H135").Value = 5569.7407
$ws.Range("J135").Value = 9964.714
$ws.Range("L135").Value = 89682.42600000001
$ws.Range("N135").Value = -94752.42600000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4146.5293
$ws.Range("I122").Value = 4239
$ws.Range("J122").Value = 4014.4285
$ws.Range("K122").Value = 12717
$ws.Range("L122").Value = 12043.2855
$ws.Range("M122").Value = -10267
$ws.Range("N122").Value = -16943.2855
$ws.Range("H132").Value = 3096.628
$ws.Range("I132").Value = 2394.9375
$ws.Range("J132").Value = 5137.909
$ws.Range("K132").Value = 7184.8125
$ws.Range("L132").Value = 15413.727
$ws.Range("M132").Value = -4654.8125
$ws.Range("N132").Value = -20473.727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4499.75
$ws.Range("I68").Value = 2559.2
$ws.Range("K68").Value = 2559.2
$ws.Range("M68").Value = -1810.2
$ws.Range("H71").Value = 4499.75
$ws.Range("I71").Value = 2559.2
$ws.Range("K71").Value = 12796
$ws.Range("M71").Value = -9052
$ws.Range("H132").Value = 5346.6294
$ws.Range("I132").Value = 2685.7368
$ws.Range("J132").Value = 11666.25
$ws.Range("K132").Value = 8057.2104
$ws.Range("L132").Value = 34998.75
$ws.Range("M132").Value = -5527.2104
$ws.Range("N132").Value = -40058.75
$ws.Range("H135").Value = 69999
$ws.Range("J135").Value = 69999
$ws.Range("L135").Value = 69999
$ws.Range("N135").Value = -80139
$ws.Range("H136").Value = 15556.667
$ws.Range("I136").Value = 2000
$ws.Range("K136").Value = 6000
$ws.Range("M136").Value = -3450

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 831
$ws.Range("J113").Value = 973.875
$ws.Range("L113").Value = 2921.625
$ws.Range("N113").Value = -7261.625
$ws.Range("H122").Value = 9378.68
$ws.Range("I122").Value = 2682
$ws.Range("J122").Value = 14640.357
$ws.Range("K122").Value = 8046
$ws.Range("L122").Value = 43921.071
$ws.Range("M122").Value = -5596
$ws.Range("N122").Value = -48821.071
$ws.Range("H132").Value = 6464.5483
$ws.Range("I132").Value = 5764.143
$ws.Range("K132").Value = 17292.429
$ws.Range("M132").Value = -14762.429
$ws.Range("H136").Value = 2967.7837
$ws.Range("I136").Value = 2372.0625
$ws.Range("J136").Value = 6780.4
$ws.Range("K136").Value = 7116.1875
$ws.Range("L136").Value = 20341.2
$ws.Range("M136").Value = -4566.1875
$ws.Range("N136").Value = -25441.2
